$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking text (e.g. "245.97") that must stay
# as text (matching the original inlineStr cell type) rather than being
# auto-converted to a number by Excel. Force the cell to Text format first.
$textForceCells = @(
    @{ Cell = 'D5'; Value = '245.97' }
    @{ Cell = 'D8'; Value = '41.96' }
    @{ Cell = 'D10'; Value = '51.01' }
    @{ Cell = 'D11'; Value = '0.0732' }
    @{ Cell = 'D12'; Value = '0.0969' }
    @{ Cell = 'D14'; Value = '12.74' }
    @{ Cell = 'D15'; Value = '0.709' }
    @{ Cell = 'D16'; Value = '4.87' }
    @{ Cell = 'D19'; Value = '72.48' }
    @{ Cell = 'D21'; Value = '243.27' }
    @{ Cell = 'D22'; Value = '12.63' }
    @{ Cell = 'D25'; Value = '2.42' }
    @{ Cell = 'D26'; Value = '2.22' }
    @{ Cell = 'D27'; Value = '164.83' }
    @{ Cell = 'D28'; Value = '8.33' }
    @{ Cell = 'D29'; Value = '18.12' }
    @{ Cell = 'D32'; Value = '1.66' }
    @{ Cell = 'D33'; Value = '4.28' }
    @{ Cell = 'D34'; Value = '0.0573' }
    @{ Cell = 'D35'; Value = '4.13' }
    @{ Cell = 'D38'; Value = '1.97' }
    @{ Cell = 'D39'; Value = '1.58' }
    @{ Cell = 'D40'; Value = '97.28' }
    @{ Cell = 'D41'; Value = '16.80' }
    @{ Cell = 'D42'; Value = '0.0658' }
    @{ Cell = 'D46'; Value = '2.31' }
    @{ Cell = 'D48'; Value = '0.0786' }
    @{ Cell = 'D50'; Value = '12.10' }
    @{ Cell = 'D51'; Value = '6.43' }
)
foreach ($item in $textForceCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
}

# Remaining cells: plain text (coin names, URLs, percentages, multi-dot
# "thousand-grouped" price strings) that Excel will not mistake for numbers.
$plainCells = @(
    @{ Cell = 'D2'; Value = '34.709.86' }
    @{ Cell = 'E2'; Value = '  -1.75%  ' }
    @{ Cell = 'D3'; Value = '1.868.96' }
    @{ Cell = 'E3'; Value = '  -2.26%  ' }
    @{ Cell = 'E4'; Value = '  -1.01%  ' }
    @{ Cell = 'E5'; Value = '  -2.75%  ' }
    @{ Cell = 'E6'; Value = '  -4.95%  ' }
    @{ Cell = 'E7'; Value = '  -1.03%  ' }
    @{ Cell = 'E8'; Value = '  +3.70%  ' }
    @{ Cell = 'E9'; Value = '  -3.84%  ' }
    @{ Cell = 'E10'; Value = '  -3.24%  ' }
    @{ Cell = 'E11'; Value = '  -0.03%  ' }
    @{ Cell = 'E12'; Value = '  -2.89%  ' }
    @{ Cell = 'D13'; Value = '2.142.75' }
    @{ Cell = 'E13'; Value = '  -2.26%  ' }
    @{ Cell = 'E14'; Value = '  +1.49%  ' }
    @{ Cell = 'E15'; Value = '  -0.65%  ' }
    @{ Cell = 'E16'; Value = '  -0.30%  ' }
    @{ Cell = 'D17'; Value = '1.879.56' }
    @{ Cell = 'E17'; Value = '  -1.95%  ' }
    @{ Cell = 'D18'; Value = '34.686.68' }
    @{ Cell = 'E18'; Value = '  -1.90%  ' }
    @{ Cell = 'E19'; Value = '  -0.63%  ' }
    @{ Cell = 'E20'; Value = '  -1.64%  ' }
    @{ Cell = 'E21'; Value = '  +0.73%  ' }
    @{ Cell = 'E22'; Value = '  -3.29%  ' }
    @{ Cell = 'E23'; Value = '  -3.40%  ' }
    @{ Cell = 'E24'; Value = '  -0.90%  ' }
    @{ Cell = 'E25'; Value = '  +4.13%  ' }
    @{ Cell = 'E26'; Value = '  -4.75%  ' }
    @{ Cell = 'E27'; Value = '  -1.76%  ' }
    @{ Cell = 'E28'; Value = '  -3.69%  ' }
    @{ Cell = 'E29'; Value = '  -3.06%  ' }
    @{ Cell = 'E30'; Value = '  -5.12%  ' }
    @{ Cell = 'D31'; Value = '4.128.56' }
    @{ Cell = 'E31'; Value = '  +0.02%  ' }
    @{ Cell = 'E32'; Value = '  +3.81%  ' }
    @{ Cell = 'E33'; Value = '  -1.42%  ' }
    @{ Cell = 'E35'; Value = '  -2.53%  ' }
    @{ Cell = 'E36'; Value = '  -1.09%  ' }
    @{ Cell = 'E37'; Value = '  -9.97%  ' }
    @{ Cell = 'B38'; Value = 'LidoDAOToken' }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' }
    @{ Cell = 'E38'; Value = '  -2.93%  ' }
    @{ Cell = 'B39'; Value = 'WEMIXToken' }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' }
    @{ Cell = 'E39'; Value = '  -20.38%  ' }
    @{ Cell = 'E40'; Value = '  -1.72%  ' }
    @{ Cell = 'E41'; Value = '  -3.82%  ' }
    @{ Cell = 'E42'; Value = '  +1.26%  ' }
    @{ Cell = 'E43'; Value = '  -0.20%  ' }
    @{ Cell = 'E44'; Value = '  -5.04%  ' }
    @{ Cell = 'D45'; Value = '1.280.08' }
    @{ Cell = 'E45'; Value = '  -4.96%  ' }
    @{ Cell = 'E46'; Value = '  -7.28%  ' }
    @{ Cell = 'E47'; Value = '  -1.02%  ' }
    @{ Cell = 'E48'; Value = '  +8.16%  ' }
    @{ Cell = 'E49'; Value = '  -2.03%  ' }
    @{ Cell = 'E50'; Value = '  +4.03%  ' }
    @{ Cell = 'E51'; Value = '  -3.60%  ' }
)
foreach ($item in $plainCells) {
    $ws.Range($item.Cell).Value = $item.Value
}

Write-Host "Applied $($textForceCells.Count) text-forced cell updates and $($plainCells.Count) plain cell updates."
